$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B92: change from text "5" to a real number 5
$ws.Range("B92").Value = 5

# Add new row 93 with the new annotation data
$ws.Range("A93").Value = "Ying Tang"
# politeness_score for row 93 stays a text "1" (not numeric), so force text format
$ws.Range("B93").NumberFormat = "@"
$ws.Range("B93").Value = "1"
$ws.Range("B93").Style = "Normal"
$ws.Range("C93").Value = "utterly unconvincing"
$ws.Range("D93").Value = "CRT"
$ws.Range("E93").Value = "THE"
$ws.Range("F93").Value = "8fd9d1eb-d55a-4b83-a989-0f77ecdd42b7"
$ws.Range("G93").Value = "r1CE9GWR-_annotated.xlsx"
$ws.Range("H93").Value = "Moreover, the discussion of supervised and unsupervised paradigms is utterly unconvincing, especially in light of the above comment on minimum-distance estimation underlying both of these paradigms."
